# Formed the consolidated report
# Update the "Absent" (column H) values so that every row's Absent flag
# is consistent with the day's attendance (H = 1 when nobody attended
# that day, H = 0 otherwise), completing the consolidated report.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H3").Value = 1
$ws.Range("H6").Value = 0
$ws.Range("H9").Value = 1
$ws.Range("H10").Value = 0
